$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table data: column A = matchup, column B = Ballgorithm pick, column C = ESPN pick
$data = @(
    @("NBA, Thursday 29th Feb 2024", "Ballgorithm", "ESPN"),
    @("Milwaukee Bucks (38-21) vs Charlotte Hornets (15-43)", "Milwaukee Bucks (77.42%)", "Milwaukee Bucks (85.4%)"),
    @("Utah Jazz (27-32) vs Orlando Magic (33-26)", "Orlando Magic (70.37%)", "Orlando Magic (68.9%)"),
    @("Atlanta Hawks (26-32) vs Brooklyn Nets (22-36)", "Brooklyn Nets (51.61%)", "Brooklyn Nets (53.9%)"),
    @("Golden State Warriors (30-27) vs New York Knicks (35-24)", "New York Knicks (66.67%)", "Golden State Warriors (53.8%)"),
    @("Oklahoma City Thunder (41-17) vs San Antonio Spurs (11-48)", "Oklahoma City Thunder (80.00%)", "Oklahoma City Thunder (84.9%)"),
    @("Houston Rockets (25-33) vs Phoenix Suns (34-24)", "Houston Rockets (66.67%)", "Phoenix Suns (75.6%)"),
    @("Miami Heat (33-25) vs Denver Nuggets (40-19)", "Denver Nuggets (82.14%)", "Denver Nuggets (62.5%)"),
    @("Washington Wizards (9-49) vs Los Angeles Lakers (33-28)", "Los Angeles Lakers (70.00%)", "Los Angeles Lakers (83.5%)")
)

# Grow the table: copy the formatting of the last existing data row (7) down
# into the two new rows (8 and 9) before writing values.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

$ws.Range("C10").Select()
